$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A94").Value = "2024-11-04 00:00:00"
$ws.Range("B94").Value = 74600
$ws.Range("C94").Value = 10461.51
$ws.Range("D94").Value = 9257.98
$ws.Range("E94").Value = 7.0956
